$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 33
$ws.Cells.Item(32, 3).Value = 11
$ws.Cells.Item(32, 4).Value = 13
$ws.Cells.Item(32, 5).Value = 37
$ws.Cells.Item(32, 6).Value = 57
$ws.Cells.Item(32, 7).Value = 94
